$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CS229Fall14")

# --- Enter the new label strings first, in the exact order the shared
#     string table records them (E3, E4, H3, H4, E6, H7..H10, H12..H15) ---
$ws.Range("E3").Value = "Mean data:"
$ws.Range("E4").Value = "Mean diffs:"
$ws.Range("H3").Value = "Upper control limit: "
$ws.Range("H4").Value = "Lower control limit:"
$ws.Range("E6").Value = "SD data:"
$ws.Range("H7").Value = "3SD UCL"
$ws.Range("H8").Value = "2SD UCL"
$ws.Range("H9").Value = "1SD UCL"
$ws.Range("H10").Value = "0.5SD UCL"
$ws.Range("H12").Value = "3SD LCL"
$ws.Range("H13").Value = "2SD LCL"
$ws.Range("H14").Value = "1SD LCL"
$ws.Range("H15").Value = "0.5SD LCL"

# --- Formulas: compute the averages / stdev first (others depend on them) ---
$ws.Range("F3").Formula = "=AVERAGE(B3:B15)"
$ws.Range("F4").Formula = "=AVERAGE(C3:C15)"
$ws.Range("F6").Formula = "=STDEV(B3:B15)"

$ws.Range("I3").Formula = "=F3+0.5*F4/1.128"
$ws.Range("I4").Formula = "=F3-3*F4/1.128"

$ws.Range("I7").Formula = '=$F$3+3*$F$6'
$ws.Range("I8").Formula = '=$F$3+2*$F$6'
$ws.Range("I9").Formula = '=$F$3+1*$F$6'
$ws.Range("I10").Formula = '=$F$3+0.5*$F$6'
$ws.Range("I12").Formula = '=$F$3-3*$F$6'
$ws.Range("I13").Formula = '=$F$3-2*$F$6'
$ws.Range("I14").Formula = '=$F$3-1*$F$6'
$ws.Range("I15").Formula = '=$F$3-0.5*$F$6'

# Column H width (matches the 17.5 stored width after Excel's internal padding)
$ws.Columns.Item(8).ColumnWidth = 16.666666666666668

# --- Make this sheet the active / selected one, with H16 selected ---
$ws.Activate()
$ws.Range("H16").Select()
